# Update column F (dSF) values per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = -2
    3  = 0
    4  = -1
    5  = 6
    6  = 3
    10 = -1
    11 = 1
    13 = 2
    16 = 2
    18 = -2
    19 = 2
    23 = 2
    24 = -2
    28 = 1
    38 = 2
    39 = 7
    42 = 7
    43 = 4
    44 = -3
    48 = -3
    49 = 6
    50 = 5
    51 = 0
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
